$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.737.70'
$ws.Range("E2").Value = '  -0.39%  '
$ws.Range("D3").Value = '2.332.87'
$ws.Range("E3").Value = '  -1.18%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.62'
$ws.Range("E5").Value = '  -1.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.659'
$ws.Range("E6").Value = '  -4.13%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '71.37'
$ws.Range("E7").Value = '  -5.95%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.580'
$ws.Range("E9").Value = '  -6.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0980'
$ws.Range("E10").Value = '  -4.21%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.76'
$ws.Range("E11").Value = '  +0.93%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '32.00'
$ws.Range("E12").Value = '  -2.43%  '
$ws.Range("E13").Value = '  -0.51%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.06'
$ws.Range("E14").Value = '  -6.22%  '
$ws.Range("D15").Value = '2.680.52'
$ws.Range("E15").Value = '  -1.29%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '16.00'
$ws.Range("E16").Value = '  -4.96%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.889'
$ws.Range("E17").Value = '  -2.91%  '
$ws.Range("D18").Value = '2.333.17'
$ws.Range("E18").Value = '  -1.48%  '
$ws.Range("D19").Value = '43.561.77'
$ws.Range("E19").Value = '  -0.74%  '
$ws.Range("E20").Value = '  -2.71%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.59'
$ws.Range("E21").Value = '  -1.11%  '
$ws.Range("B22").Value = 'Litecoin'
$ws.Range("C22").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '77.52'
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '249.97'
$ws.Range("E23").Value = '  -2.74%  '
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.87'
$ws.Range("E25").Value = '  +6.70%  '
$ws.Range("E26").Value = '  +2.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.46'
$ws.Range("E27").Value = '  -2.31%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.25'
$ws.Range("E28").Value = '  -8.30%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.26'
$ws.Range("E29").Value = '  -1.36%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '174.82'
$ws.Range("E30").Value = '  -0.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '21.98'
$ws.Range("E31").Value = '  -4.84%  '
$ws.Range("E32").Value = '  -2.47%  '
$ws.Range("E33").Value = '  -0.76%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0726'
$ws.Range("E34").Value = '  -2.99%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.02'
$ws.Range("E35").Value = '  -4.90%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.31'
$ws.Range("E36").Value = '  -0.33%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.72'
$ws.Range("E37").Value = '  -2.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.35'
$ws.Range("E38").Value = '  -3.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.31'
$ws.Range("E39").Value = '  -2.63%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.57'
$ws.Range("E40").Value = '  +23.73%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0267'
$ws.Range("E41").Value = '  -2.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '64.01'
$ws.Range("E42").Value = '  +17.65%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.16'
$ws.Range("E43").Value = '  +2.74%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.105'
$ws.Range("E44").Value = '  +4.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.72'
$ws.Range("E45").Value = '  -0.73%  '
$ws.Range("E46").Value = '  -3.93%  '
$ws.Range("E47").Value = '  -0.06%  '
$ws.Range("B48").Value = 'TrustWalletToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.21'
$ws.Range("E48").Value = '  -3.73%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.41'
$ws.Range("E49").Value = '  -3.58%  '
$ws.Range("E50").Value = '  +3.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.13'
$ws.Range("E51").Value = '  -4.48%  '
